$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (rich-text shared strings, single uniform font per string) ---
$ws.Range("A8").Value = "Volume 31   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/22/2024  Through  4/28/2024"

# --- Row 15 ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -83.333333333333
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -6.666666666666
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 180
$ws.Range("N15").Value = -30

# --- Row 16 ---
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = -6.25
$ws.Range("I16").Value = 130
$ws.Range("J16").Value = 122
$ws.Range("K16").Value = 6.55737704918
$ws.Range("L16").Value = 26.213592233009
$ws.Range("M16").Value = 64.556962025316
$ws.Range("N16").Value = -65.51724137931

# --- Row 17 ---
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 180
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = 30.434782608695
$ws.Range("I17").Value = 226
$ws.Range("J17").Value = 190
$ws.Range("K17").Value = 18.947368421052
$ws.Range("L17").Value = 14.720812182741
$ws.Range("M17").Value = 189.74358974359
$ws.Range("N17").Value = -19.858156028368

# --- Row 18 ---
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 6
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").Value = 50
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 108.333333333333
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = 36.231884057971
$ws.Range("L18").Value = -12.14953271028
$ws.Range("M18").Value = 147.368421052632
$ws.Range("N18").Value = -72.271386430678

# --- Row 19 ---
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 43.333333333333
$ws.Range("I19").Value = 181
$ws.Range("J19").Value = 145
$ws.Range("K19").Value = 24.827586206896
$ws.Range("L19").Value = 37.121212121212
$ws.Range("M19").Value = 158.571428571429
$ws.Range("N19").Value = 86.597938144329

# --- Row 20 ---
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 47
$ws.Range("H20").Value = -65.95744680851
$ws.Range("I20").Value = 78
$ws.Range("J20").Value = 153
$ws.Range("K20").Value = -49.019607843137
$ws.Range("L20").Value = -15.217391304347
$ws.Range("M20").Value = 110.810810810811
$ws.Range("N20").Value = -58.510638297872

# --- Row 21 ---
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 28.947368421052
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 175
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 724
$ws.Range("J21").Value = 698
$ws.Range("K21").Value = 3.724928366762
$ws.Range("L21").Value = 11.728395061728
$ws.Range("M21").Value = 132.051282051282
$ws.Range("N21").Value = -44.774980930587

# --- Row 22 ---
$ws.Range("M22").Value = -42.857142857142

# --- Row 23 ---
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 60
$ws.Range("F23").Value = 26
$ws.Range("H23").Value = -45.833333333333
$ws.Range("I23").Value = 127
$ws.Range("J23").Value = 142
$ws.Range("K23").Value = -10.56338028169
$ws.Range("L23").Value = 18.691588785046
$ws.Range("M23").Value = 118.965517241379

# --- Row 24 ---
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 64.705882352941
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = 32.394366197183
$ws.Range("I24").Value = 372
$ws.Range("J24").Value = 363
$ws.Range("K24").Value = 2.479338842975
$ws.Range("L24").Value = 3.910614525139
$ws.Range("M24").Value = 61.739130434782

# --- Row 25 ---
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 200
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 55
$ws.Range("J25").Value = 76
$ws.Range("K25").Value = -27.631578947368
$ws.Range("L25").Value = -48.598130841121

# --- Row 26 ---
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = -27.272727272727
$ws.Range("F26").Value = 68
$ws.Range("G26").Value = 83
$ws.Range("H26").Value = -18.072289156626
$ws.Range("I26").Value = 273
$ws.Range("J26").Value = 371
$ws.Range("K26").Value = -26.415094339622
$ws.Range("L26").Value = -14.953271028037
$ws.Range("M26").Value = 0.367647058823

# --- Row 27 ---
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 20
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = -9.090909090909
$ws.Range("L27").Value = -4.761904761904

# --- Row 28 ---
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 50
$ws.Range("I28").Value = 25
$ws.Range("J28").Value = 34
$ws.Range("K28").Value = -26.470588235294
$ws.Range("L28").Value = 13.636363636363

# --- Row 29 ---
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 6
$ws.Range("K29").Value = -40
$ws.Range("L29").Value = -40
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -75

# --- Row 30 ---
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 50
$ws.Range("I30").Value = 6
$ws.Range("K30").Value = -40
$ws.Range("L30").Value = -40
$ws.Range("M30").Value = -45.454545454545
$ws.Range("N30").Value = -75

# --- Column E bestFit width narrows to match columns C/D/F/G after the data shortened ---
$ws.Columns("E").ColumnWidth = $ws.Columns("C").ColumnWidth
